$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New labels for the 113 samples, in the same order as the existing
# row/column headers (row 1, columns B..DJ and column A, rows 2..114).
$newLabels = @("WAM63","WAM60","WAM57","WAM56","WAM55","SAM46","SAM8","SAM44","SAM41","SAM36","WAM64","SAM30","SAM29","SAM28","SAM27","SAM26","SAM25","SAM21","SAM20","SAM18","SAM17","SAM33","SAM16","SAM14","SAM12","SAM11","WAM58","ANWC92","ANWC90","ANWC88","SAM15","ANWC135","ANWC140","ANWC129","ANWC131","ANWC130","ANWC152","SAM35","ANWC159","ANWC118","ANWC119","ANWC84","ANWC94","ANWC133","ANWC136","SAM31","ANWC156","ANWC79","ANWC132","ANWC103","ANWC105","SAM32","ANWC100","SAM39","ANWC102","ANWC138","ANWC104","ANWC168","ANWC117","ANWC101","ANWC69","ANWC75","ANWC154","ANWC106","ANWC108","ANWC112","ANWC107","ANWC126","SAM45","ANWC110","ANWC85","ANWC115","ANWC137","ANWC66","ANWC116","ANWC139","SAM19","ANWC141","ANWC128","ANWC142","ANWC147","ANWC109","ANWC77","ANWC162","ANWC146","SAM34","ANWC148","ANWC150","ANWC145","ANWC151","ANWC113","ANWC153","ANWC155","SAM37","ANWC71","ANWC157","ANWC158","ANWC143","ANWC160","ANWC161","ANWC170","ANWC68","ANWC72","ANWC73","ANWC134","ANWC67","ANWC81","ANWC163","WAM59","ANWC78","SAM38","ANWC144","ANWC80")

for ($i = 0; $i -lt $newLabels.Length; $i++) {
    $label = $newLabels[$i]
    # Row 1 header (column B = index 2 .. column DJ)
    $ws.Cells.Item(1, $i + 2).Value = $label
    # Column A header (row 2 .. row 114)
    $ws.Cells.Item($i + 2, 1).Value = $label
}

# Reset the font used by the numeric matrix cells so it resolves to the
# workbook's theme body font (Calibri) instead of the hard-coded
# "Lucida Grande" font, while keeping the original font size.
$dataRange = $ws.Range("B2:DJ114")
$dataRange.Font.Name = "Calibri"
$dataRange.Font.ThemeFont = 1

# Update the active cell / selection.
[void]$ws.Range("I24").Select()
